$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 6659
$ws.Range("B2").Value = "Sr. Anthony Gabriel da Cunha"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45082
$ws.Range("G2").Value = 8981.07

# Row 3
$ws.Range("A3").Value = 68072
$ws.Range("B3").Value = "Luiza Borges"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45085
$ws.Range("G3").Value = 4943.44

# Row 4
$ws.Range("A4").Value = 23822
$ws.Range("B4").Value = "Vitória Garcia"
$ws.Range("C4").Value = "Operacoes"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 2635.71

# Row 5
$ws.Range("A5").Value = 84375
$ws.Range("B5").Value = "Ana Julia Caldeira"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45092
$ws.Range("G5").Value = 9003.700000000001

# Row 6
$ws.Range("A6").Value = 96335
$ws.Range("B6").Value = "Maria Júlia Caldeira"
$ws.Range("C6").Value = "Operacoes"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 2518.36

# Row 7
$ws.Range("A7").Value = 25876
$ws.Range("B7").Value = "Lavínia Correia"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45082
$ws.Range("G7").Value = 2696.94

# Row 8
$ws.Range("A8").Value = 81581
$ws.Range("B8").Value = "Daniel Aparecida"
$ws.Range("C8").Value = "Operacoes"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45093
$ws.Range("G8").Value = 7218.76

# Row 9
$ws.Range("A9").Value = 68168
$ws.Range("B9").Value = "Benjamin Abreu"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Consulta medica"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45093
$ws.Range("G9").Value = 9067.219999999999

# Row 10
$ws.Range("A10").Value = 29299
$ws.Range("B10").Value = "Sra. Mirella Porto"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45087
$ws.Range("G10").Value = 5211.5

# Row 11
$ws.Range("A11").Value = 25776
$ws.Range("B11").Value = "Dom Aparecida"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 4121.83
